$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-22 12:52:34"
$wsZh.Range("E3").Value = "2016-03-22 12:52:34"
$wsZh.Range("H2").Value = "2016-03-22 12:53:26"
$wsZh.Range("H3").Value = "2016-03-22 12:53:26"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-22 12:52:43"
$wsDe.Range("E3").Value = "2016-03-22 12:52:43"
$wsDe.Range("H2").Value = "2016-03-22 12:53:40"
$wsDe.Range("H3").Value = "2016-03-22 12:53:40"
